$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: surveyor is inserted at E1, notes header moves to F1
$ws.Cells.Item(1, 5).Value = "surveyor"
$ws.Cells.Item(1, 6).Value = "notes"

# Existing notes text (previously in column E) moves to column F.
$notesByRow = @{
    3  = "2 sets duplicates"
    5  = "1 tagged might be stuck to another and hidden"
    6  = "1 set duplicates"
    7  = "3 tagged missing numbers"
    8  = "1 set duplicates"
    10 = "1 empty shell"
}

foreach ($r in $notesByRow.Keys) {
    $ws.Cells.Item($r, 6).Value = $notesByRow[$r]
}

# New surveyor names for column E, data rows 2-11
$surveyors = @{
    2  = "Grace"
    3  = "Grace"
    4  = "Ava"
    5  = "Megan"
    6  = "Ava"
    7  = "Grace"
    8  = "Grace"
    9  = "Ava"
    10 = "Megan"
    11 = "Ava"
}

foreach ($r in $surveyors.Keys) {
    $ws.Cells.Item($r, 5).Value = $surveyors[$r]
}

# Update the selection to match the new state
$ws.Range("F18").Select()
